$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data cells per the latest scrape.
# Force text number format on every written cell so numeric-looking
# strings (e.g. "1.000", "0.01497") are preserved exactly as text
# instead of being coerced into floating point numbers by Excel.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.587.72'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +2.03%  '
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.662.97'
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9994'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.00%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.78'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.58%  '
# Row 6
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.01%  '
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4794'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.08%  '
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2617'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.20%  '
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06147'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.66%  '
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07074'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.40%  '
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.669.93'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.32%  '
# Row 12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +2.18%  '
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.5896'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -4.54%  '
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.374'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -4.02%  '
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '74.33'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.82%  '
# Row 16
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.01%  '
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9998'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.10%  '
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '25.567.88'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +2.07%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000006751'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +3.03%  '
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.39'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.30%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.877.33'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.49%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.423'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.22%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.644'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +2.11%  '
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.292'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.97%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '134.51'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.78%  '
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.07'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +2.48%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.407'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.51%  '
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '104.75'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +2.96%  '
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.685'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.15%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.942'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +4.09%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.648'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +2.71%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.07630'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -3.64%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9996'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.05%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04313'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -5.08%  '
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.618'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.38%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6112'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +6.00%  '
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9481'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.08%  '
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.606'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.47%  '
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.8535'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.22%  '
# Row 40
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.04%  '
# Row 41
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.01497'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -2.36%  '
# Row 42
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.873'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +3.01%  '
# Row 43
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '97.92'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.79%  '
# Row 44
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3756'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.36%  '
# Row 45
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.690'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -2.57%  '
# Row 46
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1118'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.08%  '
# Row 47
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'Aptos'
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.207'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +2.87%  '
# Row 48
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05260'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.17%  '
# Row 49
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '29.42'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.00%  '
# Row 50
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'TrueUSD'
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.002'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.13%  '
# Row 51
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'USDD'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.003'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.36%  '
